# Script2.xlsx — Nexial "expression" / "json" command catalogue update
#
# Adds the new `storeKeys(json,jsonpath,var)` command to the hidden
# "#system" sheet's `json` list (column M), and removes the redundant
# single-cell `text` column (the `text` named range already covers a
# single cell and the category name itself also lives in the `target`
# list in column A, duplicated as a one-off column Y) — collapsing the
# list of category columns so what used to be Z:AE (web, webalert,
# webcookie, ws, ws.async, xml) becomes Y:AD, and "text" is dropped from
# the `target` roster in column A as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) "json" column (M): splice in a new entry at M16 for the new
#    storeKeys(...) command, pushing the two rows below it down by one.
#    Cell-by-cell reassignment (instead of Range.Insert, which shifts the
#    *entire* row across every column) keeps the rest of the sheet intact.
# ---------------------------------------------------------------------
$ws.Range("M18").Value = $ws.Range("M17").Value()
$ws.Range("M17").Value = $ws.Range("M16").Value()
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 2) "target" column (A): drop the "text" entry (was A25), shifting the
#    remaining entries (web, webalert, webcookie, ws, ws.async, xml) up
#    by one row and clearing the now-vacant last row (A31).
# ---------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $ws.Range("A" + $r).Value = $ws.Range("A" + ($r + 1)).Value()
}
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------
# 3) Drop column Y (the single-cell "text" list) entirely, shifting
#    columns Z:AE (web, webalert, webcookie, ws, ws.async, xml) one
#    column to the left, to Y:AD.
# ---------------------------------------------------------------------
$ws.Range("Y1").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 4) Update the defined names so each still points at its (now shifted)
#    data range.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
